$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "CMS COMPONENTS" paragraph: mention the new MySQL database component,
#    fix the old MYSQL -> MySQL casing, and add a trailing sentence
#    describing the database.
# ---------------------------------------------------------------------------

# Insert ", and a MySQL database" right after "...two NodeJS applications"
# (before the period that follows it).
$rng = $d.Content
$rng.Find.Execute("two NodeJS applications.", $false, $false, $false, $false, $false, $true, 1, $false, `
    "two NodeJS applications, and a MySQL database.", 2) | Out-Null

# Fix the old ALL-CAPS "MYSQL" reference to the properly-cased "MySQL".
$rng = $d.Content
$rng.Find.Execute("MYSQL database for display", $false, $false, $false, $false, $false, $true, 1, $false, `
    "MySQL database for display", 2) | Out-Null

# Append a new closing sentence about the MySQL database to the end of the
# paragraph (right after "...to be viewed.").
$rng = $d.Content
$rng.Find.Execute("to be viewed.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("  The MySQL database is where content to be displayed on both NodeJS applications is stored.")

# ---------------------------------------------------------------------------
# 2. Turn the plain-text front end URL into a real hyperlink.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("http://localhost:3001/", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($rng, "http://localhost:3001/") | Out-Null
$rng.Font.Size = 12

# ---------------------------------------------------------------------------
# 3. Remove the old trailing blank paragraph, then add the new "MySQL
#    Database Information" section in its place, finishing with the
#    _GoBack bookmark on the very last line.
# ---------------------------------------------------------------------------

# Drop the old bookmark (it currently sits right after the URL).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Merge away the empty paragraph that used to trail the document.
$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1)
$last = $d.Paragraphs.Last
$mergeRng = $d.Range($secondLast.Range.End - 1, $last.Range.End)
$mergeRng.Delete()

$lines = @(
    "MySQL Database Information:",
    "To sign into the database, please see the following credentials:",
    "Website: phpmyadmin.co",
    "Server: sql9.freesqldatabase.com",
    "Username: sql9229224",
    "Password: 6m2d4QZdzj"
)

$cur = $d.Paragraphs.Last.Range
$cur.Collapse(0)

foreach ($line in $lines) {
    $pos = $cur.End
    $cur.InsertParagraphAfter()
    $cur = $d.Range($pos + 1, $pos + 1)
    $cur.InsertAfter($line)
}

# Re-create the _GoBack bookmark at the very end of the document (right
# after the password line, which is now the last line of content).
$endPos = $d.Paragraphs.Last.Range.End - 1
$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
